# Update the build version timestamp throughout the workbook
# Old: February 03 2026 17.29.55 EST
# New: February 03 2026 18.05.36 EST

$wb = $excel.ActiveWorkbook

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

# "About" sheet updates
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: $newVersion"

$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Wangpo Coal Mine, China, M1217, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)`""

# "Boundaries and methane sources" sheet updates (build_version column, rows 2-13)
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 13; $row++) {
    $dataSheet.Cells.Item($row, 19).Value = $newVersion
}
